$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted at row 93 ("Fruta / hortaliza,
# semanal"). Insert a blank row there which pushes the existing rows 93-144
# down to 94-145 (this also naturally extends the used range / dimension to
# R145, matching row 144's old data landing intact on the new row 145).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly record. The
# non-varying columns (market/region/category/quality/unit/origin/etc.) are
# identical to every other row in this block.
$ws.Range("A93").Value = 4
$ws.Range("B93").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C93").Value = "Los Lagos"
$ws.Range("D93").Value = 44460
$ws.Range("E93").Value = 10
$ws.Range("F93").Value = 100112043
$ws.Range("G93").Value = "Pepino ensalada"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 160
$ws.Range("K93").Value = 17500
$ws.Range("L93").Value = 17500
$ws.Range("M93").Value = 17500
$ws.Range("N93").Value = "$/caja 60 unidades"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 292
$ws.Range("Q93").Value = 60
$ws.Range("R93").Value = "Hortaliza"
